$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 467, shifting existing rows 467-509 down to 468-510.
$ws.Rows.Item(467).Insert()

# Populate the newly inserted row 467 with the new data record.
$ws.Cells.Item(467, 1).Value = 6
$ws.Cells.Item(467, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(467, 3).Value = "Metropolitana"
$ws.Cells.Item(467, 4).Value = 44858
$ws.Cells.Item(467, 5).Value = 13
$ws.Cells.Item(467, 6).Value = 100112043
$ws.Cells.Item(467, 7).Value = "Pepino ensalada"
$ws.Cells.Item(467, 8).Value = "Sin especificar"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 430
$ws.Cells.Item(467, 11).Value = 16000
$ws.Cells.Item(467, 12).Value = 18000
$ws.Cells.Item(467, 13).Value = 17000
$ws.Cells.Item(467, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(467, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(467, 16).Value = 283
$ws.Cells.Item(467, 17).Value = 60
$ws.Cells.Item(467, 18).Value = "Hortaliza"
